$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7603300213813782
$ws.Range("B1").Value = 1.164723992347717
$ws.Range("C1").Value = 2.272536754608154
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.81972599029541
